# The workbook has two sheets: "Exportación" (first) and "Importación" (second).
# On the "Exportación" sheet, two pairs of rows had their data corrected/re-sorted:
#  - the row for Pedimento 4001081 (08/05/2024, Clave R1) and the row for
#    Pedimento 4000834 (07/05/2024, Clave RT) were swapped between rows 3 and 4
#    (columns C "Pedimento", E "Clave" and F "Fecha Pago")
#  - the row for Pedimento 4001130 (Clave R1) and the row for Pedimento 4001070
#    (Clave RT) were swapped between rows 11 and 12 (columns C "Pedimento" and
#    E "Clave"; both rows already shared the same "Fecha Pago" so F is untouched)
#
# The swaps are performed with Range.Copy(Destination) (instead of assigning
# Range.Value directly) so that date-like text such as "08/05/2024" is carried
# over as-is and is not re-interpreted/auto-converted into an Excel date serial
# number, which is what would happen if a literal date-like string were typed
# (or assigned through .Value) into the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exportación")

# --- Swap rows 3 and 4 (columns C, E, F) using row 200 as scratch space ---
$ws.Range("C3").Copy($ws.Range("C200"))
$ws.Range("E3").Copy($ws.Range("E200"))
$ws.Range("F3").Copy($ws.Range("F200"))

$ws.Range("C4").Copy($ws.Range("C3"))
$ws.Range("E4").Copy($ws.Range("E3"))
$ws.Range("F4").Copy($ws.Range("F3"))

$ws.Range("C200").Copy($ws.Range("C4"))
$ws.Range("E200").Copy($ws.Range("E4"))
$ws.Range("F200").Copy($ws.Range("F4"))

# --- Swap rows 11 and 12 (columns C, E) using row 200 as scratch space ---
$ws.Range("C11").Copy($ws.Range("C200"))
$ws.Range("E11").Copy($ws.Range("E200"))

$ws.Range("C12").Copy($ws.Range("C11"))
$ws.Range("E12").Copy($ws.Range("E11"))

$ws.Range("C200").Copy($ws.Range("C12"))
$ws.Range("E200").Copy($ws.Range("E12"))

# Clean up the scratch cells
$ws.Range("C200:F200").Clear()
